$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 267, shifting existing rows 267-318 down to 268-319
$ws.Rows(267).Insert()

# Populate the newly inserted row 267 with the new record
$ws.Cells.Item(267, 1).Value = 10
$ws.Cells.Item(267, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(267, 3).Value = 'La Araucanía'
$ws.Cells.Item(267, 4).Value = 44637
$ws.Cells.Item(267, 5).Value = 9
$ws.Cells.Item(267, 6).Value = 'Fruta'
$ws.Cells.Item(267, 7).Value = 100108
$ws.Cells.Item(267, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(267, 9).Value = 100108002
$ws.Cells.Item(267, 10).Value = 'Mango'
$ws.Cells.Item(267, 11).Value = 'Sin especificar'
$ws.Cells.Item(267, 12).Value = 'Primera'
$ws.Cells.Item(267, 13).Value = 310
$ws.Cells.Item(267, 14).Value = 7500
$ws.Cells.Item(267, 15).Value = 7500
$ws.Cells.Item(267, 16).Value = 7500
$ws.Cells.Item(267, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(267, 18).Value = 'Perú'
$ws.Cells.Item(267, 19).Value = 1875
$ws.Cells.Item(267, 20).Value = 4
